# Updates odds values in Sheet1, matching the commit "Atualizando o arquivo XLSX".
# Only numeric odds cells in rows 2, 4, 5 and 6 change; everything else stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (El Gaish vs El Ismaily)
$ws.Range("G2").Value  = 2.35
$ws.Range("H2").Value  = 2.77
$ws.Range("I2").Value  = 3.4
$ws.Range("J2").Value  = 3.05
$ws.Range("L2").Value  = 4.1
$ws.Range("W2").Value  = 2.12
$ws.Range("Y2").Value  = 5.7
$ws.Range("AA2").Value = 9.75
$ws.Range("AB2").Value = 25
$ws.Range("AC2").Value = 24
$ws.Range("AF2").Value = 5.6
$ws.Range("AG2").Value = 18
$ws.Range("AK2").Value = 16
$ws.Range("AL2").Value = 12.5
$ws.Range("AM2").Value = 50
$ws.Range("AO2").Value = 55

# Row 4 (Al Fateh vs Al Raed)
$ws.Range("H4").Value  = 4
$ws.Range("I4").Value  = 4.75
$ws.Range("K4").Value  = 2.3
$ws.Range("M4").Value  = 1.03
$ws.Range("N4").Value  = 10
$ws.Range("O4").Value  = 1.2
$ws.Range("P4").Value  = 4.33
$ws.Range("Q4").Value  = 1.67
$ws.Range("R4").Value  = 2.15
$ws.Range("S4").Value  = 2.63
$ws.Range("T4").Value  = 1.44
$ws.Range("U4").Value  = 1.33
$ws.Range("V4").Value  = 3.25
$ws.Range("W4").Value  = 1.73
$ws.Range("X4").Value  = 2
$ws.Range("Y4").Value  = 8
$ws.Range("Z4").Value  = 8.5
$ws.Range("AB4").Value = 13
$ws.Range("AD4").Value = 23
$ws.Range("AE4").Value = 13
$ws.Range("AF4").Value = 8
$ws.Range("AH4").Value = 41
$ws.Range("AL4").Value = 15
$ws.Range("AN4").Value = 34

# Row 5 (Al Ittihad vs Al Riyadh)
$ws.Range("O5").Value  = 1.11
$ws.Range("P5").Value  = 6
$ws.Range("S5").Value  = 1.91
$ws.Range("T5").Value  = 1.8
$ws.Range("AD5").Value = 19
$ws.Range("AE5").Value = 23

# Row 6 (Al Shabab vs Al Orubah)
$ws.Range("G6").Value  = 1.2
$ws.Range("H6").Value  = 7
$ws.Range("K6").Value  = 2.88
$ws.Range("M6").Value  = 1.02
$ws.Range("N6").Value  = 12
$ws.Range("O6").Value  = 1.11
$ws.Range("P6").Value  = 6
$ws.Range("Q6").Value  = 1.4
$ws.Range("R6").Value  = 2.75
$ws.Range("S6").Value  = 1.91
$ws.Range("T6").Value  = 1.8
$ws.Range("Z6").Value  = 7.5
$ws.Range("AE6").Value = 21
$ws.Range("AJ6").Value = 29
